$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.542.14'
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").Value = '2.597.24'
$ws.Range("E3").Value = '  +10.15%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.32'
$ws.Range("E5").Value = '  +1.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.00'
$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.600'
$ws.Range("E7").Value = '  +5.23%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  +13.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.54'
$ws.Range("E10").Value = '  +11.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("E11").Value = '  +4.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.18'
$ws.Range("E12").Value = '  +14.49%  '

$ws.Range("D13").Value = '2.992.70'
$ws.Range("E13").Value = '  +10.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +1.90%  '

$ws.Range("D15").Value = '2.588.77'
$ws.Range("E15").Value = '  +9.72%  '

$ws.Range("E16").Value = '  +10.91%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.90'
$ws.Range("E17").Value = '  +9.41%  '

$ws.Range("D18").Value = '46.666.83'
$ws.Range("E18").Value = '  +1.24%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  +4.52%  '

$ws.Range("E20").Value = '  +4.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.67'
$ws.Range("E21").Value = '  +9.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.18'
$ws.Range("E22").Value = '  +5.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.21'
$ws.Range("E23").Value = '  +3.92%  '

$ws.Range("E24").Value = '  +4.86%  '

$ws.Range("E25").Value = '  +13.71%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.12'
$ws.Range("E26").Value = '  +33.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.49'
$ws.Range("E28").Value = '  +6.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.30'
$ws.Range("E29").Value = '  +4.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '39.57'
$ws.Range("E30").Value = '  -0.52%  '

$ws.Range("B31").Value = 'LidoDAOToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.73'
$ws.Range("E31").Value = '  -0.61%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.12'
$ws.Range("E32").Value = '  +10.41%  '

$ws.Range("E33").Value = '  +22.57%  '

$ws.Range("E34").Value = '  +5.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0834'
$ws.Range("E35").Value = '  +7.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '149.89'
$ws.Range("E36").Value = '  +2.48%  '

$ws.Range("E37").Value = '  +3.85%  '

$ws.Range("E38").Value = '  +4.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.20'
$ws.Range("E39").Value = '  +6.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '15.74'
$ws.Range("E40").Value = '  +4.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.63'
$ws.Range("E41").Value = '  +12.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0323'
$ws.Range("E42").Value = '  +7.48%  '

$ws.Range("D43").Value = '2.028.79'
$ws.Range("E43").Value = '  +6.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.48'
$ws.Range("E44").Value = '  +27.24%  '

$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.74'
$ws.Range("E46").Value = '  -0.69%  '

$ws.Range("E47").Value = '  -1.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '109.04'
$ws.Range("E48").Value = '  +11.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.08'
$ws.Range("E49").Value = '  +10.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.201'
$ws.Range("E50").Value = '  +7.54%  '

$ws.Range("D51").Value = '2.849.30'
$ws.Range("E51").Value = '  +10.06%  '
